$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the "Curso (semestre ideal)" line to drop "EQD (8), "
# ------------------------------------------------------------------
$d.Content.Find.Execute("Curso (semestre ideal): EQD (8), EQN (10)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Curso (semestre ideal): EQN (10)", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Reorder the "Requisitos" bullet list: the LOQ4057 line (currently
#    first) must move to become the last line (after LOQ4002). The
#    three lines live as three separate runs (text + manual line
#    break) inside a single paragraph. We cut the first run (text +
#    its break) and re-insert it at the end of the paragraph, right
#    before the paragraph mark - this keeps it as its own run, same
#    as in the target document.
# ------------------------------------------------------------------

# Locate the paragraph that contains the requirement list via the
# unique "LOQ4057" token.
$findRng = $d.Content
$findRng.Find.Execute("LOQ4057", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hitStart = $findRng.Start

$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($hitStart -ge $p.Range.Start -and $hitStart -lt $p.Range.End) {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start
$full = $target.Range.Text

# Manual line breaks show up as chr(11) in Range.Text; the first one
# marks the end of the LOQ4057 run (including the break itself).
$breakOffsets = @()
for ($i = 0; $i -lt $full.Length; $i++) {
    if ([int][char]$full[$i] -eq 11) {
        $breakOffsets += $i
    }
}

$run1 = $d.Range($pStart, $pStart + $breakOffsets[0] + 1)
$run1Text = $run1.Text
$run1.Cut() | Out-Null

# Insert right before the paragraph mark (end of the now-shrunk
# paragraph), after the LOQ4002 run/break.
$newPEnd = $target.Range.End
$insertPoint = $d.Range($newPEnd - 1, $newPEnd - 1)
$insertPoint.InsertAfter($run1Text) | Out-Null
